# Apply updated odds values to the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.55
$ws.Range("G2").Value = 1.56
$ws.Range("H2").Value = 6.4
$ws.Range("I2").Value = 6.6
$ws.Range("J2").Value = 4.8
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 1.37
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 2.26
$ws.Range("Q2").Value = 1.77
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 1.87
$ws.Range("U2").Value = 2.12
$ws.Range("V2").Value = 1.17
$ws.Range("W2").Value = 2.78
$ws.Range("Z2").Value = 65
$ws.Range("AA2").Value = 200
$ws.Range("AB2").Value = 8.800000000000001
$ws.Range("AC2").Value = 10.5
$ws.Range("AD2").Value = 23
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 9.199999999999999
$ws.Range("AG2").Value = 9.4
$ws.Range("AH2").Value = 22
$ws.Range("AI2").Value = 80
$ws.Range("AJ2").Value = 14
$ws.Range("AK2").Value = 15
$ws.Range("AL2").Value = 30
$ws.Range("AM2").Value = 110
$ws.Range("AN2").Value = 7.4
$ws.Range("AO2").Value = 130

# Row 3
$ws.Range("C3").Value = "09:30:00"
$ws.Range("G3").Value = 2.02
$ws.Range("J3").Value = 1.09
$ws.Range("M3").Value = 1.04
$ws.Range("W3").Value = 1.99

# Row 4
$ws.Range("F4").Value = 1.35
$ws.Range("G4").Value = 1.4
$ws.Range("H4").Value = 14.5
$ws.Range("I4").Value = 17.5
$ws.Range("K4").Value = 5.1
$ws.Range("L4").Value = 1.51
$ws.Range("N4").Value = 2.98
$ws.Range("P4").Value = 1.67
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.24
$ws.Range("T4").Value = 2.8
$ws.Range("W4").Value = 3.5
$ws.Range("X4").Value = 14
$ws.Range("Y4").Value = 32
$ws.Range("AB4").Value = 5.7
$ws.Range("AC4").Value = 12.5
$ws.Range("AD4").Value = 1000
$ws.Range("AF4").Value = 6.4
$ws.Range("AH4").Value = 55
$ws.Range("AJ4").Value = 11
$ws.Range("AL4").Value = 190
$ws.Range("AN4").Value = 10

# Row 5
$ws.Range("F5").Value = 2.14
$ws.Range("G5").Value = 2.6
$ws.Range("H5").Value = 2.98
$ws.Range("I5").Value = 3.75
$ws.Range("J5").Value = 3.25
$ws.Range("K5").Value = 4.1
$ws.Range("L5").Value = 1.36
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 3.3
$ws.Range("P5").Value = 1.89
$ws.Range("Q5").Value = 1.8
$ws.Range("R5").Value = 1.37
$ws.Range("S5").Value = 2.88
$ws.Range("T5").Value = 1.66
$ws.Range("U5").Value = 2.1
$ws.Range("V5").Value = 1.38
$ws.Range("W5").Value = 1.62
$ws.Range("X5").Value = 16.5
$ws.Range("Y5").Value = 14.5
$ws.Range("Z5").Value = 25
$ws.Range("AA5").Value = 60
$ws.Range("AB5").Value = 12
$ws.Range("AC5").Value = 9
$ws.Range("AD5").Value = 15
$ws.Range("AE5").Value = 40
$ws.Range("AF5").Value = 17.5
$ws.Range("AG5").Value = 12.5
$ws.Range("AH5").Value = 18.5
$ws.Range("AI5").Value = 48
$ws.Range("AJ5").Value = 36
$ws.Range("AK5").Value = 27
$ws.Range("AL5").Value = 40
$ws.Range("AM5").Value = 200
$ws.Range("AN5").Value = 20
$ws.Range("AO5").Value = 38

# Row 6
$ws.Range("G6").Value = 8.6
$ws.Range("I6").Value = 1.5
$ws.Range("N6").Value = 3.9
$ws.Range("O6").Value = 1.32
$ws.Range("P6").Value = 2
$ws.Range("Q6").Value = 1.96
$ws.Range("T6").Value = 2.14
$ws.Range("U6").Value = 1.81
$ws.Range("X6").Value = 16
$ws.Range("Z6").Value = 8
$ws.Range("AH6").Value = 28
$ws.Range("AI6").Value = 38
$ws.Range("AJ6").Value = 300
$ws.Range("AO6").Value = 8.199999999999999

# Row 8
$ws.Range("F8").Value = 1.93
$ws.Range("G8").Value = 2.18
$ws.Range("H8").Value = 3.9
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 3.3
$ws.Range("K8").Value = 4
$ws.Range("L8").Value = 1.42
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 3.4
$ws.Range("O8").Value = 1.34
$ws.Range("P8").Value = 1.81
$ws.Range("Q8").Value = 1.99
$ws.Range("R8").Value = 1.3
$ws.Range("S8").Value = 3.55
$ws.Range("T8").Value = 1.81
$ws.Range("U8").Value = 1.94
$ws.Range("V8").Value = 1.28
$ws.Range("W8").Value = 1.84
$ws.Range("X8").Value = 16
$ws.Range("Y8").Value = 17
$ws.Range("Z8").Value = 40
$ws.Range("AA8").Value = 120
$ws.Range("AB8").Value = 9.199999999999999
$ws.Range("AC8").Value = 8.6
$ws.Range("AD8").Value = 22
$ws.Range("AE8").Value = 65
$ws.Range("AF8").Value = 14.5
$ws.Range("AG8").Value = 12.5
$ws.Range("AH8").Value = 23
$ws.Range("AJ8").Value = 29
$ws.Range("AK8").Value = 27
$ws.Range("AL8").Value = 44
$ws.Range("AN8").Value = 21
$ws.Range("AO8").Value = 85
